# Begin portfolio optimization program
# - Select G20 on the QSTK sheet (moves active-cell selection there; it
#   also loses tabSelected once a new sheet becomes active below).
# - Add a new worksheet "HW1 - AssessPorfolio" after the last sheet (QSTK),
#   which becomes the active sheet/tab.
# - Populate the new sheet with the homework notes content.

$wb = $excel.ActiveWorkbook

# Update stale selection left on QSTK before it stops being the active tab.
$qstk = $wb.Worksheets.Item("QSTK")
[void]$qstk.Range("G20").Select()

# Add the new sheet at the very end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "HW1 - AssessPorfolio"

# Content
$newSheet.Range("A1").Value = "4 equities"
$newSheet.Range("A2").Value = "Determine the best combo to get the highest Sharpe Ratio"
$newSheet.Range("A3").Value = "Hold all year - 2011"

$newSheet.Range("A5").Value = "Tutorial1.py"
$newSheet.Range("A6").Value = "Tutorial3.py - do not use method in tuturial3 for HW1"
$newSheet.Range("B7").Value = ">> no rebalancing"

$newSheet.Range("A9").Value = "PART 1"
$newSheet.Range("A11").Value = "PART 2"

$newSheet.Range("A13").Value = "PART 3"
$newSheet.Range("B13").Value = 1
$newSheet.Range("C13").Value = "use simulate() function to optimize portfolio based on allocations"
$newSheet.Range("B14").Value = 2
$newSheet.Range("C14").Value = "created nested for loops to run different scenarios"

# Narrow the "numbering" column (B) like the other note sheets in this
# workbook. The engine quantizes ColumnWidth to whole display pixels, so
# this lands on the closest reachable width to the source file's 2.6640625.
$newSheet.Columns.Item(2).ColumnWidth = 1.75

[void]$newSheet.Range("C15").Select()
